# Generate Report for Handback
# Updates timestamps / priority values produced by a fresh handback report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview.Range("G2").Value = "2016-08-24 18:23:30"
$wsOverview.Range("G5").Value = "2016-08-24 18:23:30"

# --- zh-cn sheet ---
# Priority column (E): ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

# Correspond Handoff Datetime column (H)
$wsZhCn.Range("H2").Value = "2016-08-24 18:23:25"
$wsZhCn.Range("H5").Value = "2016-08-24 18:23:25"

# Correspond Handback DateTime column (K)
$wsZhCn.Range("K2").Value = "2016-08-24 18:23:42"
$wsZhCn.Range("K5").Value = "2016-08-24 18:23:42"

# --- de-de sheet ---
# Priority column (E): ht -> mt (same shared string as zh-cn's Priority column)
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# Correspond Handoff Datetime column (H) - shares the string with Overview!G
$wsDeDe.Range("H2").Value = "2016-08-24 18:23:30"
$wsDeDe.Range("H5").Value = "2016-08-24 18:23:30"

# Correspond Handback DateTime column (K)
$wsDeDe.Range("K2").Value = "2016-08-24 18:23:49"
$wsDeDe.Range("K5").Value = "2016-08-24 18:23:49"
